$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.049.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.908.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "370.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.67%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.87%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0836"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.372.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.918.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.948"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.095.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0945"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.96%  "
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("E26").Value = "  +4.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.172"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.103"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.91%  "
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0422"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.67%  "
$ws.Range("E42").Value = "  -7.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.113"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "117.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("E46").Value = "  -3.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.045.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.208.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.240"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.13%  "
